$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 80: correct the date/time value in column A ---
$ws.Range("A80").Value = 45457.2916666667

# --- Row 81: new data row appended after row 80 ---
$ws.Range("A81").Value = 45460.296412037
# Match A80's date/time number format & style (copy formats only)
$ws.Range("A80").Copy()
$ws.Range("A81").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("B81").Value = 1500
$ws.Range("C81").Value = 2.98000001907349
$ws.Range("D81").Value = 2.98000001907349
$ws.Range("E81").Value = 2.98000001907349
$ws.Range("F81").Value = 2.98000001907349

# G81 stores the adj_close as text (matches existing column G convention)
$ws.Range("G81").NumberFormat = "@"
$ws.Range("G81").Value = "2.98000001907349"
$ws.Range("G81").ClearFormats()

$ws.Range("H81").Value = "ESPE.MI"
